# feat. narration skip & quiz check
#
# Rewrites the NarrationData sheet: the 4 placeholder test rows are
# replaced with the full intro / tutorial / forest narration script
# (19 entries), column C is widened to fit the longer Korean text, and
# the sheet's used range grows from A1:E6 to A1:E21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C needs to be much wider to hold the long narration sentences.
$ws.Columns(3).ColumnWidth = 92

# Row 3-4: intro narration
$ws.Range("B3").Value  = 'intro_01'
$ws.Range("C3").Value  = 'whitecane에 오신 여러분 환영합니다.'
$ws.Range("B4").Value  = 'intro_02'
$ws.Range("C4").Value  = '지금부터 여러분의 시각 일부를 제한하고 청각과 촉감에 의존하여 떠나는 색다른 모험으로 안내하겠습니다.'

# Row 5-14: tutorial narration
$ws.Range("B5").Value  = 'tutorial01'
$ws.Range("C5").Value  = '시작하기 앞서 게임의 방법을 소개해드리겠습니다'
$ws.Range("B6").Value  = 'tutorial02'
$ws.Range("C6").Value  = '물체를 맞추는 방법에 대해 알려드리겠습니다'
$ws.Range("B7").Value  = 'tutorial03'
$ws.Range("C7").Value  = 'a를 꾹 눌러 녹음을 진행하고 버튼을 땝니다'
$ws.Range("B8").Value  = 'tutorial04'
$ws.Range("C8").Value  = '한번해볼까요?'
$ws.Range("B9").Value  = 'tutorial05'
$ws.Range("C9").Value  = '잘하셨습니다'
$ws.Range("B10").Value = 'tutorial06'
$ws.Range("C10").Value = '이제 지팡이를 드리겠습니다'
$ws.Range("B11").Value = 'tutorial07'
$ws.Range("C11").Value = '지팡이가 닿는 곳은 부분적으로 밝아집니다.'
$ws.Range("B12").Value = 'tutorial08'
$ws.Range("C12").Value = '앞에 무엇이 잇는지 몰라 망설여질 때 유용하게 활용해보세요'
$ws.Range("B13").Value = 'tutorial09'
$ws.Range("C13").Value = '소리가 들리는 곳으로 발걸음을 옮겨주시기를 바랍니다. 즐거운 여정이 되시길 바랍니다.'
$ws.Range("B14").Value = 'tutorial10'
$ws.Range("C14").Value = '즐거운 여정이 되시길 바랍니다.'

# Row 15-21: forest (foreset) narration
$ws.Range("B15").Value = 'foreset_01'
$ws.Range("C15").Value = '여러분은 현재 드넓은 숲 한가운데에 서 있습니다.'
$ws.Range("B16").Value = 'foreset_02'
$ws.Range("C16").Value = '이곳에서 나가기 위해서는 폭포를 찾아가야합니다.'
$ws.Range("B17").Value = 'foreset_03'
$ws.Range("C17").Value = '폭포를 찾아가는 길에는 예상치 못한 퀴즈들이 기다리고 있습니다.'
$ws.Range("B18").Value = 'foreset_04'
$ws.Range("C18").Value = '소리를 듣고 소리의 정체가 무엇인지 맞추면 되는 간단한 퀴즈입니다.'
$ws.Range("B19").Value = 'foreset_05'
$ws.Range("C19").Value = '답을 맞추시면 폭포 위치에 대한 힌트를 얻을 수 있습니다.'
$ws.Range("B20").Value = 'foreset_06'
$ws.Range("C20").Value = '이제 폭포소리를 따라 발걸음을 옮겨주시기를 바랍니다.'
$ws.Range("B21").Value = 'foreset_07'
$ws.Range("C21").Value = '즐거운 여정이 되시길 바랍니다.'
